$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this weekly block (row 390), shifting the
# existing rows 390-413 down to 392-415.
$ws.Rows("390:391").Insert()

# New row 390 - Ají / Inferno / Primera, week of 2023-04-25 (serial 45041)
$ws.Range("A390").Value = 8
$ws.Range("B390").Value = "Terminal La Palmera de La Serena"
$ws.Range("C390").Value = "Coquimbo"
$ws.Range("D390").Value = 45041
$ws.Range("E390").Value = 4
$ws.Range("F390").Value = 100112021
$ws.Range("G390").Value = "Ají"
$ws.Range("H390").Value = "Inferno"
$ws.Range("I390").Value = "Primera"
$ws.Range("J390").Value = 460
$ws.Range("K390").Value = 11000
$ws.Range("L390").Value = 12000
$ws.Range("M390").Value = 11500
$ws.Range("N390").Value = "$/caja 15 kilos"
$ws.Range("O390").Value = "Provincia de Limarí"
$ws.Range("P390").Value = 767
$ws.Range("Q390").Value = 15
$ws.Range("R390").Value = "Hortaliza"

# New row 391 - Ají / Inferno / Segunda, same week (serial 45041)
$ws.Range("A391").Value = 8
$ws.Range("B391").Value = "Terminal La Palmera de La Serena"
$ws.Range("C391").Value = "Coquimbo"
$ws.Range("D391").Value = 45041
$ws.Range("E391").Value = 4
$ws.Range("F391").Value = 100112021
$ws.Range("G391").Value = "Ají"
$ws.Range("H391").Value = "Inferno"
$ws.Range("I391").Value = "Segunda"
$ws.Range("J391").Value = 300
$ws.Range("K391").Value = 7000
$ws.Range("L391").Value = 8000
$ws.Range("M391").Value = 7500
$ws.Range("N391").Value = "$/caja 15 kilos"
$ws.Range("O391").Value = "Provincia de Limarí"
$ws.Range("P391").Value = 500
$ws.Range("Q391").Value = 15
$ws.Range("R391").Value = "Hortaliza"
